# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap the order of "Santa Lucia" and "Timor Oriental" in row 202/203 ---
# Row 202 currently shows "Santa Lucia", row 203 shows "Timor Oriental".
# Target: row 202 -> "Timor Oriental", row 203 -> "Santa Lucia".
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# --- Update "last updated" timestamp string (row 1 title) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Agosto de 2020 a las 12:11"

# --- Update per-country statistics ---
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5874295
$ws.Range("C4").Value = 149
$ws.Range("D4").Value = 3167164
$ws.Range("E4").Value = 2526526
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 180605

# Row 14: Iran
$ws.Range("B14").Value = 361150
$ws.Range("C14").Value = 2245
$ws.Range("D14").Value = 311365
$ws.Range("E14").Value = 29009
$ws.Range("G14").Value = 133
$ws.Range("H14").Value = 20776

# Row 18: Banglades
$ws.Range("B18").Value = 297083
$ws.Range("C18").Value = 2485
$ws.Range("D18").Value = 182875
$ws.Range("E18").Value = 110225
$ws.Range("G18").Value = 42
$ws.Range("H18").Value = 3983

# Row 42: Rumania
$ws.Range("B42").Value = 79330
$ws.Range("C42").Value = 825
$ws.Range("D42").Value = 35517
$ws.Range("E42").Value = 40504
$ws.Range("G42").Value = 37
$ws.Range("H42").Value = 3309

# Row 89: Libia
$ws.Range("B89").Value = 11009
$ws.Range("C89").Value = 572
$ws.Range("D89").Value = 1096
$ws.Range("E89").Value = 9714
$ws.Range("G89").Value = 11
$ws.Range("H89").Value = 199

# Row 93: Guinea
$ws.Range("B93").Value = 9013
$ws.Range("C93").Value = 46
$ws.Range("D93").Value = 7823
$ws.Range("E93").Value = 1136
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 54

# Row 101: Finlandia
$ws.Range("B101").Value = 7938
$ws.Range("C101").Value = 18
$ws.Range("E101").Value = 504

# Row 126: Sri Lanka
$ws.Range("D126").Value = 2811
$ws.Range("E126").Value = 130

# Row 172: Birmania
$ws.Range("B172").Value = 463
$ws.Range("C172").Value = 13
$ws.Range("E172").Value = 116

# Row 194: Liechtenstein
$ws.Range("B194").Value = 100
$ws.Range("C194").Value = 1
$ws.Range("D194").Value = 90
$ws.Range("E194").Value = 9
